$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 148.59
$ws.Range("I15").Value = 148.59
$ws.Range("K15").Value = 445.77
$ws.Range("M15").Value = -276.77
$ws.Range("H29").Value = 1266.9231
$ws.Range("I29").Value = 636.2857
$ws.Range("J29").Value = 2002.6666
$ws.Range("K29").Value = 1908.8571
$ws.Range("L29").Value = 6007.9998
$ws.Range("M29").Value = -1627.8571
$ws.Range("N29").Value = -6569.9998
$ws.Range("H86").Value = 2083.0
$ws.Range("I86").Value = 699.53845
$ws.Range("J86").Value = 5680.0
$ws.Range("K86").Value = 699.53845
$ws.Range("L86").Value = 5680.0
$ws.Range("M86").Value = 423.46155
$ws.Range("N86").Value = -7926.0
$ws.Range("H89").Value = 2083.0
$ws.Range("I89").Value = 699.53845
$ws.Range("J89").Value = 5680.0
$ws.Range("K89").Value = 3497.69225
$ws.Range("L89").Value = 28400.0
$ws.Range("M89").Value = 2118.30775
$ws.Range("N89").Value = -39632.0
$ws.Range("H96").Value = 635.7143
$ws.Range("I96").Value = 300.0
$ws.Range("K96").Value = 900.0
$ws.Range("M96").Value = 473.0
$ws.Range("H106").Value = 2943.111
$ws.Range("I106").Value = 2943.111
$ws.Range("K106").Value = 2943.111
$ws.Range("M106").Value = -2312.111
$ws.Range("H112").Value = 1411.84
$ws.Range("J112").Value = 1421.4949
$ws.Range("L112").Value = 4264.4847
$ws.Range("N112").Value = -6480.4847
$ws.Range("H116").Value = 1117388.9
$ws.Range("J116").Value = 10100.0
$ws.Range("L116").Value = 10100.0
$ws.Range("N116").Value = -16984.0
$ws.Range("H132").Value = 27891346.0
$ws.Range("I132").Value = 33468534.0
$ws.Range("J132").Value = 5410.0
$ws.Range("K132").Value = 100405602.0
$ws.Range("L132").Value = 16230.0
$ws.Range("M132").Value = -100403072.0
$ws.Range("N132").Value = -21290.0
$ws.Range("H137").Value = 3415.5918
$ws.Range("I137").Value = 2840.425
$ws.Range("J137").Value = 5971.8887
$ws.Range("K137").Value = 8521.275000000001
$ws.Range("L137").Value = 17915.6661
$ws.Range("M137").Value = -5971.275000000001
$ws.Range("N137").Value = -23015.6661
$ws.Range("H138").Value = 2477.34
$ws.Range("I138").Value = 1267.5312
$ws.Range("J138").Value = 3046.6619
$ws.Range("K138").Value = 3802.5936
$ws.Range("L138").Value = 9139.985700000001
$ws.Range("M138").Value = 1337.4064
$ws.Range("N138").Value = -19419.9857

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9639.444
$ws.Range("I32").Value = 6672.8945
$ws.Range("J32").Value = 14763.484
$ws.Range("K32").Value = 6672.8945
$ws.Range("L32").Value = 14763.484
$ws.Range("M32").Value = -6385.8945
$ws.Range("N32").Value = -15337.484
$ws.Range("H110").Value = 1000.0
$ws.Range("I110").Value = 1000.0
$ws.Range("K110").Value = 1000.0
$ws.Range("M110").Value = 1045.0
$ws.Range("H122").Value = 2705.6553
$ws.Range("I122").Value = 1466.8889
$ws.Range("J122").Value = 4732.727
$ws.Range("K122").Value = 4400.6667
$ws.Range("L122").Value = 14198.181
$ws.Range("M122").Value = -1950.6667
$ws.Range("N122").Value = -19098.181
$ws.Range("H132").Value = 2032.8108
$ws.Range("I132").Value = 983.9655
$ws.Range("J132").Value = 5834.875
$ws.Range("K132").Value = 2951.8965
$ws.Range("L132").Value = 17504.625
$ws.Range("M132").Value = -421.8964999999998
$ws.Range("N132").Value = -22564.625
$ws.Range("H137").Value = 0.0
$ws.Range("J137").Value = 0.0
$ws.Range("L137").Value = 0.0
$ws.Range("N137").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 103.5
$ws.Range("I7").Value = 105.0
$ws.Range("J7").Value = 90.0
$ws.Range("K7").Value = 105.0
$ws.Range("L7").Value = 90.0
$ws.Range("M7").Value = 8.0
$ws.Range("N7").Value = -316.0
$ws.Range("H31").Value = 3834.946
$ws.Range("I31").Value = 1494.7
$ws.Range("J31").Value = 6588.1763
$ws.Range("K31").Value = 1494.7
$ws.Range("L31").Value = 6588.1763
$ws.Range("M31").Value = -1199.7
$ws.Range("N31").Value = -7178.1763
$ws.Range("H34").Value = 3834.946
$ws.Range("I34").Value = 1494.7
$ws.Range("J34").Value = 6588.1763
$ws.Range("K34").Value = 1494.7
$ws.Range("L34").Value = 6588.1763
$ws.Range("M34").Value = -1292.7
$ws.Range("N34").Value = -6992.1763
$ws.Range("H94").Value = 2889.125
$ws.Range("I94").Value = 2700.0
$ws.Range("J94").Value = 2952.1667
$ws.Range("K94").Value = 2700.0
$ws.Range("L94").Value = 2952.1667
$ws.Range("M94").Value = -2249.0
$ws.Range("N94").Value = -3854.1667
$ws.Range("H132").Value = 5389.875
$ws.Range("I132").Value = 4887.273
$ws.Range("K132").Value = 14661.819
$ws.Range("M132").Value = -12131.819

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2644.75
$ws.Range("J5").Value = 3979.5
$ws.Range("L5").Value = 11938.5
$ws.Range("N5").Value = -12162.5
$ws.Range("H26").Value = 7941.9473
$ws.Range("I26").Value = 13522.0
$ws.Range("J26").Value = 2919.9
$ws.Range("K26").Value = 40566.0
$ws.Range("L26").Value = 8759.7
$ws.Range("M26").Value = -40278.0
$ws.Range("N26").Value = -9335.7
$ws.Range("H113").Value = 755.75
$ws.Range("J113").Value = 986.3333
$ws.Range("L113").Value = 2958.9999
$ws.Range("N113").Value = -7298.9999
$ws.Range("H135").Value = 2644.75
$ws.Range("J135").Value = 3979.5
$ws.Range("L135").Value = 35815.5
$ws.Range("N135").Value = -40885.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 29110.889
$ws.Range("J4").Value = 29110.889
$ws.Range("L4").Value = 29110.889
$ws.Range("N4").Value = -29334.889
$ws.Range("H102").Value = 2816.56
$ws.Range("I102").Value = 2322.348
$ws.Range("J102").Value = 8500.0
$ws.Range("K102").Value = 2322.348
$ws.Range("L102").Value = 8500.0
$ws.Range("M102").Value = -700.348
$ws.Range("N102").Value = -11744.0

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6619.8066
$ws.Range("I132").Value = 2428.0
$ws.Range("J132").Value = 7625.84
$ws.Range("K132").Value = 7284.0
$ws.Range("L132").Value = 22877.52
$ws.Range("M132").Value = -4754.0
$ws.Range("N132").Value = -27937.52

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 46777.723
$ws.Range("I4").Value = 100000.0
$ws.Range("J4").Value = 4199.9
$ws.Range("K4").Value = 100000.0
$ws.Range("L4").Value = 4199.9
$ws.Range("M4").Value = -99887.0
$ws.Range("N4").Value = -4425.9
$ws.Range("H136").Value = 2631.2563
$ws.Range("I136").Value = 758.96295
$ws.Range("J136").Value = 6843.9165
$ws.Range("K136").Value = 2276.88885
$ws.Range("L136").Value = 20531.7495
$ws.Range("M136").Value = 273.1111500000002
$ws.Range("N136").Value = -25631.7495
